# Updated cryptos list with fresh price / 1h-volume-change figures.
# Only column D (Price) and column E (Volume(1h)) values change;
# values are written as text so cell typing matches the source feed
# (Price/Volume columns hold formatted strings, not numeric values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.999.03'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '3.107.34'
$ws.Range("D5").Value = "'578.47"
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").Value = "'172.91"
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  -0.56%  '
$ws.Range("E9").Value = '  +0.49%  '
$ws.Range("D10").Value = "'0.153"
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("D13").Value = "'36.80"
$ws.Range("E14").Value = '  -1.68%  '
$ws.Range("D15").Value = '3.621.90'
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("D16").Value = '66.964.08'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("D18").Value = '3.108.94'
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("D19").Value = "'16.42"
$ws.Range("E19").Value = '  +0.87%  '
$ws.Range("D20").Value = "'491.16"
$ws.Range("E20").Value = '  +2.10%  '
$ws.Range("D21").Value = "'0.704"
$ws.Range("E21").Value = '  -1.56%  '
$ws.Range("D22").Value = "'7.88"
$ws.Range("E22").Value = '  +4.41%  '
$ws.Range("D23").Value = "'83.87"
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").Value = "'13.10"
$ws.Range("E24").Value = '  -1.21%  '
$ws.Range("E25").Value = '  -3.47%  '
$ws.Range("D26").Value = "'10.46"
$ws.Range("E26").Value = '  +4.06%  '
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").Value = "'7.90"
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("E29").Value = '  -1.42%  '
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("D31").Value = "'28.31"
$ws.Range("E31").Value = '  -1.80%  '
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("E33").Value = '  -6.10%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").Value = "'5.84"
$ws.Range("E35").Value = '  -0.86%  '
$ws.Range("E36").Value = '  -1.86%  '
$ws.Range("D37").Value = "'47.17"
$ws.Range("E38").Value = '  -3.87%  '
$ws.Range("E39").Value = '  -2.92%  '
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("E41").Value = '  -2.55%  '
$ws.Range("D42").Value = "'384.09"
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = '2.806.24'
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("E44").Value = '  -7.95%  '
$ws.Range("E45").Value = '  -2.63%  '
$ws.Range("D46").Value = "'135.51"
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").Value = "'24.80"
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("E49").Value = '  -1.90%  '
$ws.Range("E50").Value = '  -1.30%  '
$ws.Range("D51").Value = "'6.71"
$ws.Range("E51").Value = '  -2.08%  '
